$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: B1 header text changes from "in_alguma_coisa" to "in_Saldo_Da_Conta"
$ws.Range("B1").Value = "in_Saldo_Da_Conta"

# Row 2: B2 gets a new text value "1500" (quote-prefixed, kept as text)
$ws.Range("B2").Value = "'1500"

# New, currently-empty placeholder cells C2/D2 that share B2's (quote-prefix) style
$ws.Range("C2:D2").Value = "'"
$ws.Range("C2:D2").ClearContents()

# Column sizing to fit the new header / calculator contents
$ws.Columns("B").AutoFit()
$ws.Columns("D").AutoFit()

# Move the active selection to C2, matching where editing continues next
$ws.Range("C2").Select() | Out-Null
